# This document contains a single 20x5 table of simple arithmetic
# problems (e.g. "65-28=37"). Each cell's text is replaced with its
# new value, addressed by (row, column) so that cells sharing
# identical original text (e.g. the two "73-14=59" cells) are updated
# independently and correctly. Direct Range.Text assignment is used
# instead of Find/Replace because this runtime's Find.Execute searches
# and replaces across the whole document rather than being scoped to
# the Range it was invoked on.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "81-76=5"
$t.Cell(1,2).Range.Text = "39+59=98"
$t.Cell(1,3).Range.Text = "13-5=8"
$t.Cell(1,4).Range.Text = "19+17=36"
$t.Cell(1,5).Range.Text = "81-4=77"
$t.Cell(2,1).Range.Text = "52-15=37"
$t.Cell(2,2).Range.Text = "22-6=16"
$t.Cell(2,3).Range.Text = "80-49=31"
$t.Cell(2,4).Range.Text = "83-48=35"
$t.Cell(2,5).Range.Text = "15+79=94"
$t.Cell(3,1).Range.Text = "27+34=61"
$t.Cell(3,2).Range.Text = "51-24=27"
$t.Cell(3,3).Range.Text = "47+5=52"
$t.Cell(3,4).Range.Text = "68+6=74"
$t.Cell(3,5).Range.Text = "5+38=43"
$t.Cell(4,1).Range.Text = "79+16=95"
$t.Cell(4,2).Range.Text = "90-36=54"
$t.Cell(4,3).Range.Text = "42-4=38"
$t.Cell(4,4).Range.Text = "15+28=43"
$t.Cell(4,5).Range.Text = "29+19=48"
$t.Cell(5,1).Range.Text = "7+54=61"
$t.Cell(5,2).Range.Text = "90-41=49"
$t.Cell(5,3).Range.Text = "57-49=8"
$t.Cell(5,4).Range.Text = "61-58=3"
$t.Cell(5,5).Range.Text = "28+23=51"
$t.Cell(6,1).Range.Text = "26+66=92"
$t.Cell(6,2).Range.Text = "31-14=17"
$t.Cell(6,3).Range.Text = "28+17=45"
$t.Cell(6,4).Range.Text = "80-43=37"
$t.Cell(6,5).Range.Text = "24-16=8"
$t.Cell(7,1).Range.Text = "90-87=3"
$t.Cell(7,2).Range.Text = "36+57=93"
$t.Cell(7,3).Range.Text = "64-19=45"
$t.Cell(7,4).Range.Text = "73-69=4"
$t.Cell(7,5).Range.Text = "39+22=61"
$t.Cell(8,1).Range.Text = "37+49=86"
$t.Cell(8,2).Range.Text = "68+28=96"
$t.Cell(8,3).Range.Text = "56-17=39"
$t.Cell(8,4).Range.Text = "19+69=88"
$t.Cell(8,5).Range.Text = "70-17=53"
$t.Cell(9,1).Range.Text = "66+7=73"
$t.Cell(9,2).Range.Text = "61-46=15"
$t.Cell(9,3).Range.Text = "91-45=46"
$t.Cell(9,4).Range.Text = "6+49=55"
$t.Cell(9,5).Range.Text = "73-27=46"
$t.Cell(10,1).Range.Text = "70-48=22"
$t.Cell(10,2).Range.Text = "90-47=43"
$t.Cell(10,3).Range.Text = "65-56=9"
$t.Cell(10,4).Range.Text = "19+9=28"
$t.Cell(10,5).Range.Text = "67-49=18"
$t.Cell(11,1).Range.Text = "6+67=73"
$t.Cell(11,2).Range.Text = "39+26=65"
$t.Cell(11,3).Range.Text = "71-46=25"
$t.Cell(11,4).Range.Text = "90-89=1"
$t.Cell(11,5).Range.Text = "68-9=59"
$t.Cell(12,1).Range.Text = "44-25=19"
$t.Cell(12,2).Range.Text = "74-38=36"
$t.Cell(12,3).Range.Text = "26+17=43"
$t.Cell(12,4).Range.Text = "19+52=71"
$t.Cell(12,5).Range.Text = "29+37=66"
$t.Cell(13,1).Range.Text = "91-88=3"
$t.Cell(13,2).Range.Text = "57+4=61"
$t.Cell(13,3).Range.Text = "55-49=6"
$t.Cell(13,4).Range.Text = "57+25=82"
$t.Cell(13,5).Range.Text = "13-7=6"
$t.Cell(14,1).Range.Text = "71-43=28"
$t.Cell(14,2).Range.Text = "15+39=54"
$t.Cell(14,3).Range.Text = "15+46=61"
$t.Cell(14,4).Range.Text = "58+5=63"
$t.Cell(14,5).Range.Text = "72-37=35"
$t.Cell(15,1).Range.Text = "70-39=31"
$t.Cell(15,2).Range.Text = "17+38=55"
$t.Cell(15,3).Range.Text = "84-67=17"
$t.Cell(15,4).Range.Text = "78+17=95"
$t.Cell(15,5).Range.Text = "5+48=53"
$t.Cell(16,1).Range.Text = "45+6=51"
$t.Cell(16,2).Range.Text = "46-17=29"
$t.Cell(16,3).Range.Text = "23-7=16"
$t.Cell(16,4).Range.Text = "29+65=94"
$t.Cell(16,5).Range.Text = "90-23=67"
$t.Cell(17,1).Range.Text = "58+26=84"
$t.Cell(17,2).Range.Text = "43-14=29"
$t.Cell(17,3).Range.Text = "20-8=12"
$t.Cell(17,4).Range.Text = "62+19=81"
$t.Cell(17,5).Range.Text = "57+37=94"
$t.Cell(18,1).Range.Text = "73-7=66"
$t.Cell(18,2).Range.Text = "67-48=19"
$t.Cell(18,3).Range.Text = "94-38=56"
$t.Cell(18,4).Range.Text = "45+18=63"
$t.Cell(18,5).Range.Text = "90-14=76"
$t.Cell(19,1).Range.Text = "47+18=65"
$t.Cell(19,2).Range.Text = "76-29=47"
$t.Cell(19,3).Range.Text = "79+19=98"
$t.Cell(19,4).Range.Text = "15+28=43"
$t.Cell(19,5).Range.Text = "15+76=91"
$t.Cell(20,1).Range.Text = "94-39=55"
$t.Cell(20,2).Range.Text = "70-17=53"
$t.Cell(20,3).Range.Text = "75-68=7"
$t.Cell(20,4).Range.Text = "72-67=5"
$t.Cell(20,5).Range.Text = "36+47=83"
